$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Sales Scenario" - update raw input cells; dependent formulas recalc
# ---------------------------------------------------------------------------
$wsSales = $wb.Worksheets.Item("Sales Scenario")

# Monthly cost-per-customer inputs (rows 21-30, columns C:F). Column G/H unchanged.
$wsSales.Range("C21").Value = 15000000
$wsSales.Range("D21").Value = 12000000
$wsSales.Range("E21").Value = 11456000
$wsSales.Range("F21").Value = 10000000

$wsSales.Range("C22").Value = 16500000
$wsSales.Range("D22").Value = 13000000
$wsSales.Range("E22").Value = 12345000
$wsSales.Range("F22").Value = 12000000

$wsSales.Range("C23").Value = 17500000
$wsSales.Range("D23").Value = 15438000
$wsSales.Range("E23").Value = 13456000
$wsSales.Range("F23").Value = 12657000

$wsSales.Range("C24").Value = 18500000
$wsSales.Range("D24").Value = 16754000
$wsSales.Range("E24").Value = 14500000
$wsSales.Range("F24").Value = 13567000

$wsSales.Range("C25").Value = 19560000
$wsSales.Range("D25").Value = 17456000
$wsSales.Range("E25").Value = 15500000
$wsSales.Range("F25").Value = 14567000

$wsSales.Range("C26").Value = 21000000
$wsSales.Range("D26").Value = 18767900
$wsSales.Range("E26").Value = 16500000
$wsSales.Range("F26").Value = 15456000

$wsSales.Range("C27").Value = 21657000
$wsSales.Range("D27").Value = 19456000
$wsSales.Range("E27").Value = 17500000
$wsSales.Range("F27").Value = 16456000

$wsSales.Range("C28").Value = 23546000
$wsSales.Range("D28").Value = 21000000
$wsSales.Range("E28").Value = 18500000
$wsSales.Range("F28").Value = 17657800

$wsSales.Range("C29").Value = 25987000
$wsSales.Range("D29").Value = 22000000
$wsSales.Range("E29").Value = 19500000
$wsSales.Range("F29").Value = 18678000

$wsSales.Range("C30").Value = 26567000
$wsSales.Range("D30").Value = 22345000
$wsSales.Range("E30").Value = 21000000
$wsSales.Range("F30").Value = 19789000

# Column C (rows 22-30) picks up the "no top border" look already used by C21
$wsSales.Range("C21").Copy()
$wsSales.Range("C22:C30").PasteSpecial(-4122)  # xlPasteFormats

# Sunk-cost / investment figures (rows 35-44, column D)
$wsSales.Range("D35").Value = 750000000
$wsSales.Range("D36").Value = 976000000
$wsSales.Range("D37").Value = 1567000000
$wsSales.Range("D38").Value = 1998000000
$wsSales.Range("D39").Value = 2300000000
$wsSales.Range("D40").Value = 2567000000
$wsSales.Range("D41").Value = 2879000000
$wsSales.Range("D42").Value = 3240000000
$wsSales.Range("D43").Value = 3456000000
$wsSales.Range("D44").Value = 3765000000
